# LOT2026.xlsx update (2022-09-26 build):
# The course-syllabus export dropped one row (the old row 13, a stray
# B/C-only row that duplicated the professor's name with no label in
# column A) causing every subsequent row to shift up by one. Several
# cells below that point also picked up new/shuffled text content.
#
# Step 1: remove the stray row 13 so everything below shifts up,
# which reproduces the row-height / row-index layout seen in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(13).Delete()

# Step 2: fix up the handful of cells whose text content changed/shuffled
# as part of the same edit (beyond the pure row shift).
$ws.Range("B10").Value = "3403572 - Ismael Maciel de Mancilha"
$ws.Range("C10").Value = "3403572 - Ismael Maciel de Mancilha"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# (copy, not Value=, so the "01/01/2018" text isn't re-parsed as a date)
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

$ws.Range("B18").Value = "3403572 - Ismael Maciel de Mancilha"
$ws.Range("C18").Value = "3403572 - Ismael Maciel de Mancilha"

$ws.Range("B19").Value = "A avaliação será feita por meio de provas escritas."
$ws.Range("C19").Value = "A avaliação será feita por meio de provas escritas."

$ws.Range("B20").Value = "2 provas (P1 + P2), sendo que a NF = (P1 + P2) / 2"
$ws.Range("C20").Value = "2 provas (P1 + P2), sendo que a NF = (P1 + P2) / 2"

$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
